# Adapt column header formatting to respective input file names.
# "_old" -> "_FV2410", "_new" -> "_FV2504" (column K "diff" stays the same),
# turn the data range into an Excel Table, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row shared strings in place ------------------
$oldHeaders = @("Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old","Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old")
$newHeaders = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeaders[$i]
}

$oldHeaders2 = @("Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new","Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new")
$newHeaders2 = @("Segmentname_FV2504","Segmentgruppe_FV2504","Segment_FV2504","Datenelement_FV2504","Segment ID_FV2504","Code_FV2504","Qualifier_FV2504","Beschreibung_FV2504","Bedingungsausdruck_FV2504","Bedingung_FV2504")

for ($i = 0; $i -lt $oldHeaders2.Count; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $newHeaders2[$i]
}

# --- 2. Turn A1:U63 into an Excel Table (Table1) ------------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U63"), $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row (pane split after row 1) ------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
